$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to text format so numeric-looking strings
# (e.g. "0.9991", "310.82") are stored as text, matching the source data
# which uses inline/shared strings rather than numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.027.91'
$ws.Range("E2").Value = '  -2.80%  '

$ws.Range("D3").Value = '1.728.26'
$ws.Range("E3").Value = '  -1.60%  '

$ws.Range("D4").Value = '0.9991'
$ws.Range("E4").Value = '  -0.52%  '

$ws.Range("D5").Value = '310.82'
$ws.Range("E5").Value = '  -5.41%  '

$ws.Range("D6").Value = '0.9993'
$ws.Range("E6").Value = '  -0.39%  '

$ws.Range("D7").Value = '0.4871'
$ws.Range("E7").Value = '  +7.65%  '

$ws.Range("D8").Value = '0.3516'
$ws.Range("E8").Value = '  +0.75%  '

$ws.Range("D9").Value = '42.19'
$ws.Range("E9").Value = '  +0.42%  '

$ws.Range("D10").Value = '0.07278'
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").Value = '1.052'
$ws.Range("E11").Value = '  -3.36%  '

$ws.Range("D12").Value = '0.9984'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").Value = '19.97'
$ws.Range("E13").Value = '  -2.96%  '

$ws.Range("D14").Value = '5.886'
$ws.Range("E14").Value = '  -1.48%  '

$ws.Range("D15").Value = '1.719.29'
$ws.Range("E15").Value = '  -2.31%  '

$ws.Range("D16").Value = '6.877'
$ws.Range("E16").Value = '  -3.94%  '

$ws.Range("D17").Value = '87.13'
$ws.Range("E17").Value = '  -5.09%  '

$ws.Range("D18").Value = '0.00001040'
$ws.Range("E18").Value = '  -1.05%  '

$ws.Range("D19").Value = '0.06405'
$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").Value = '0.9996'
$ws.Range("E20").Value = '  -0.29%  '

$ws.Range("D21").Value = '16.57'
$ws.Range("E21").Value = '  -1.74%  '

$ws.Range("D22").Value = '5.673'
$ws.Range("E22").Value = '  -1.12%  '

$ws.Range("D23").Value = '27.059.71'
$ws.Range("E23").Value = '  -2.84%  '

$ws.Range("D24").Value = '10.85'
$ws.Range("E24").Value = '  -2.91%  '

$ws.Range("D25").Value = '2.081'
$ws.Range("E25").Value = '  -3.32%  '

$ws.Range("D26").Value = '154.28'
$ws.Range("E26").Value = '  -4.38%  '

$ws.Range("D27").Value = '20.04'
$ws.Range("E27").Value = '  +0.12%  '

$ws.Range("D28").Value = '1.917.96'
$ws.Range("E28").Value = '  -2.33%  '

$ws.Range("D29").Value = '2.082'
$ws.Range("E29").Value = '  -3.06%  '

$ws.Range("D30").Value = '121.29'
$ws.Range("E30").Value = '  -1.36%  '

$ws.Range("D31").Value = '1.040'
$ws.Range("E31").Value = '  -2.84%  '

$ws.Range("D32").Value = '0.09331'
$ws.Range("E32").Value = '  +0.90%  '

$ws.Range("D33").Value = '3.595'
$ws.Range("E33").Value = '  -1.25%  '

$ws.Range("D34").Value = '5.395'
$ws.Range("E34").Value = '  -2.63%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '0.02195'
$ws.Range("E35").Value = '  -2.82%  '

$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").Value = '0.05921'
$ws.Range("E36").Value = '  -2.73%  '

$ws.Range("D37").Value = '1.436'
$ws.Range("E37").Value = '  +4.24%  '

$ws.Range("D38").Value = '11.04'
$ws.Range("E38").Value = '  -5.74%  '

$ws.Range("D39").Value = '0.2004'
$ws.Range("E39").Value = '  -3.25%  '

$ws.Range("D40").Value = '4.779'
$ws.Range("E40").Value = '  -2.15%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '0.6017'
$ws.Range("E41").Value = '  -2.74%  '

$ws.Range("B42").Value = 'Frax'
$ws.Range("C42").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D42").Value = '0.9995'
$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("D43").Value = '1.095'
$ws.Range("E43").Value = '  -7.16%  '

$ws.Range("D44").Value = '7.553'
$ws.Range("E44").Value = '  -2.90%  '

$ws.Range("D45").Value = '12.76'
$ws.Range("E45").Value = '  -2.16%  '

$ws.Range("D46").Value = '3.585'
$ws.Range("E46").Value = '  -3.89%  '

$ws.Range("D47").Value = '0.5645'
$ws.Range("E47").Value = '  -2.50%  '

$ws.Range("D48").Value = '118.51'
$ws.Range("E48").Value = '  -2.86%  '

$ws.Range("D49").Value = '1.846'
$ws.Range("E49").Value = '  -4.11%  '

$ws.Range("D50").Value = '1.113'
$ws.Range("E50").Value = '  -0.77%  '

$ws.Range("D51").Value = '0.06657'
$ws.Range("E51").Value = '  -1.93%  '

# Restore default styling (removes the temporary text number format
# added above) so cells keep their original style index.
$ws.Range("D2:E51").Style = "Normal"
